$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four new label cells (new shared strings) --------------------
$ws.Range("I2").Value  = "wzgledne"
$ws.Range("I6").Value  = "bezwzgledne"
$ws.Range("I10").Value = "mieszane wiersze"
$ws.Range("I14").Value = "mieszane kolumny"

# --- Move the "Duplo" / LEFT() demo from J11:K11 to M2:N2 ------------------
$ws.Range("M2").Value   = "Duplo"
$ws.Range("N2").Formula = "=LEFT(M2)"
$ws.Range("J11").Clear()
$ws.Range("K11").Clear()

# --- Move the VAT-rate cell from J14 to M5, and its dependant formula ------
# --- from L14 to O5 ---------------------------------------------------------
$ws.Range("M5").Value        = 0.23
$ws.Range("M5").NumberFormat = "0%"
$ws.Range("O5").Formula      = "=100*st_vat"
$ws.Range("J14").Clear()
$ws.Range("L14").Clear()

# Re-point the defined name st_vat at its new home (M5). Deleting and
# re-adding (rather than just editing RefersTo) makes sure every formula
# that uses the name gets re-bound to the new cell.
$wb.Names.Item("st_vat").Delete()
$wb.Names.Add("st_vat", "=Arkusz1!`$M`$5")

# --- Match the saved selection shown in the workbook -----------------------
$ws.Range("I15").Select() | Out-Null
